$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps 2 RGB")

# --- Insert first block of 3 new rows (before old row 208 "V13.07.15") ---
$ws.Rows("208:210").Insert()

$ws.Range("A208").Value = "V13.07.6"
$ws.Range("B208").Value = 16.4045
$ws.Range("C208").Value = 0.8559
$ws.Range("D208").Value = 18.1467
$ws.Range("E208").Value = 0.8634
$ws.Range("F208").Value = 27.9334
$ws.Range("G208").Value = 0.8419
$ws.Range("H208").Value = 24.1549
$ws.Range("I208").Value = 0.9198
$ws.Range("J208").Value = 0.4786
$ws.Range("L208").Value = 15.253
$ws.Range("M208").Value = 0.7151

$ws.Range("A209").Value = "V13.07.7"
$ws.Range("B209").Value = 15.1744
$ws.Range("C209").Value = 0.8474
$ws.Range("D209").Value = 17.5357
$ws.Range("E209").Value = 0.8476
$ws.Range("F209").Value = 26.2995
$ws.Range("G209").Value = 0.7953
$ws.Range("H209").Value = 26.4851
$ws.Range("I209").Value = 0.9277
$ws.Range("J209").Value = 0.4159
$ws.Range("L209").Value = 14.9364
$ws.Range("M209").Value = 0.7113

$ws.Range("A210").Value = "V13.07.8"
$ws.Range("B210").Value = 15.6777
$ws.Range("C210").Value = 0.8442
$ws.Range("D210").Value = 18.5737
$ws.Range("E210").Value = 0.8837
$ws.Range("F210").Value = 26.4048
$ws.Range("G210").Value = 0.81136
$ws.Range("H210").Value = 26.15963
$ws.Range("I210").Value = 0.9418
$ws.Range("J210").Value = 0.4283
$ws.Range("L210").Value = 15.8205
$ws.Range("M210").Value = 0.7193

# --- Insert second block of 3 new rows (before "V13.10.15", now at row 213) ---
$ws.Rows("213:215").Insert()

$ws.Range("A213").Value = "V13.10.6"
$ws.Range("B213").Value = 17.202
$ws.Range("C213").Value = 0.8111
$ws.Range("D213").Value = 18.8538
$ws.Range("E213").Value = 0.848
$ws.Range("F213").Value = 35.4712
$ws.Range("G213").Value = 0.9326
$ws.Range("H213").Value = 23.576
$ws.Range("I213").Value = 0.8262
$ws.Range("J213").Value = 0.4828
$ws.Range("L213").Value = 12.6548
$ws.Range("M213").Value = 0.6995

$ws.Range("A214").Value = "V13.10.7"
$ws.Range("B214").Value = 9.6016
$ws.Range("C214").Value = 0.7267
$ws.Range("D214").Value = 15.2173
$ws.Range("E214").Value = 0.7831
$ws.Range("F214").Value = 30.9236
$ws.Range("G214").Value = 0.9023
$ws.Range("H214").Value = 23.0795
$ws.Range("I214").Value = 0.797
$ws.Range("J214").Value = 0.3983
$ws.Range("L214").Value = 14.5869
$ws.Range("M214").Value = 0.7231

$ws.Range("A215").Value = "V13.10.8"
$ws.Range("B215").Value = 9.6027
$ws.Range("C215").Value = 0.7233
$ws.Range("D215").Value = 15.4608
$ws.Range("E215").Value = 0.7831
$ws.Range("F215").Value = 30.564
$ws.Range("G215").Value = 0.8871
$ws.Range("H215").Value = 22.5272
$ws.Range("I215").Value = 0.7894
$ws.Range("J215").Value = 0.3905
$ws.Range("L215").Value = 14.7304
$ws.Range("M215").Value = 0.7231

# --- Update sheet view (pane/selection) ---
$ws.Application.ActiveWindow.ScrollRow = 198
$ws.Range("L208").Select()

# --- Update workbook window size/position (maximized) ---
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Left = -120
$excel.ActiveWindow.Width = 38640
$excel.ActiveWindow.Height = 21120
